$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Split the first paragraph ("<- Back to Home" hyperlink, style
#    FirstParagraph) into two paragraphs:
#      - a new paragraph with a "Home" hyperlink (style FirstParagraph)
#      - the original "<- Back to Home" hyperlink, now styled BodyText
# -----------------------------------------------------------------

$first = $d.Paragraphs(1)
$firstRange = $first.Range

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:hyperlink r:id="rId9"><w:r><w:rPr><w:color w:val="0066CC"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Home</w:t></w:r></w:hyperlink></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:hyperlink r:id="rId9"><w:r><w:rPr><w:color w:val="0066CC"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">' + [char]0x2190 + ' Back to Home</w:t></w:r></w:hyperlink></w:p>'

$firstRange.InsertXML($newXml)

# InsertXML does not preserve the character style reference (rStyle),
# only the direct formatting, so re-apply the "Hyperlink" character
# style to each of the two runs we just created (run-only range, so
# the paragraph style is left untouched).
$p1 = $d.Paragraphs(1)
$run1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$run1.Style = "Hyperlink"

$p2 = $d.Paragraphs(2)
$run2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$run2.Style = "Hyperlink"

# -----------------------------------------------------------------
# 2. Update every table's preferred width.
#    All tables switch from 100% (pct/5000) to automatic (auto/0) and
#    have their tblStyle/tblW order normalized, except the table that
#    already has an explicit fixed tblLayout, which keeps pct/5000
#    (its tblPr element order is simply normalized the same way).
# -----------------------------------------------------------------

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables($i)

    if ($i -eq 8) {
        # Table with explicit <w:tblLayout w:type="fixed"/>: keep 100% width.
        $t.PreferredWidthType = 2
        $t.PreferredWidth = 250
    } else {
        $t.PreferredWidthType = 1
        $t.PreferredWidth = 0
    }
}

Write-Output "edit complete"
